$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iAU_TC_ID_169"
$ws.Range("B2").Value = "@RegressionA Validation of Exam Print "
